# Commit: "Update [User Class] seuccess !!"
# Mark a batch of rooms as "Occupied" (column I = Status) on the
# RoomStock sheet — rows 20, 22, 23, 25, 29, 30, 31, 36, 41 flip from
# "Unoccupied" to "Occupied".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToMarkOccupied = @(20, 22, 23, 25, 29, 30, 31, 36, 41)
$statusCol = 9  # column I = "Status"

foreach ($r in $rowsToMarkOccupied) {
    $ws.Cells.Item($r, $statusCol).Value = "Occupied"
}
